$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Remove the 5th student row (Van Thuong) entirely - brings the table down to 3 students
$ws1.Rows.Item(5).Delete()

# Row 2 - was "Pink Rose", now "Vladimir Putin"
$ws1.Range("B2").Value = "Vladimir"
$ws1.Range("C2").Value = "Putin"
$ws1.Range("D2").Value = 36244.42661986111
$ws1.Range("E2").Value = "Male"
$ws1.Range("F2").Value = "12589"
$ws1.Range("G2").Value = "Moscow"
$ws1.Range("H2").Value = "D:\Tai lieu mon hoc 2024\Lập trình trực quan\Putin.jpeg"

# Row 3 - was "Tap Can Binh", now "Ma Tieu Dao"
$ws1.Range("B3").Value = "Ma"
$ws1.Range("C3").Value = "Tieu Dao"
$ws1.Range("D3").Value = 31237.429842002315
$ws1.Range("E3").Value = "Female"
$ws1.Range("F3").Value = "2020"
$ws1.Range("G3").Value = "Quang Chau"
$ws1.Range("H3").Value = "D:\Tai lieu mon hoc 2024\Lập trình trực quan\Ma Tieu Dao.jpeg"

# Row 4 - was "Ma Tieu Dao", now "Bi Bi Dong"
$ws1.Range("B4").Value = "Bi Bi"
$ws1.Range("C4").Value = "Dong"
$ws1.Range("D4").Value = -16092.430743344907
$ws1.Range("E4").Value = "Female"
$ws1.Range("F4").Value = "1966"
$ws1.Range("G4").Value = "Quang Nam"
$ws1.Range("H4").Value = "D:\Tai lieu mon hoc 2024\Lập trình trực quan\Bi Bi Dong.jpeg"

# Sheet2 summary stats - male ratio now 0%, female ratio now 100%
$ws2.Range("B2").Value = 0
$ws2.Range("C2").Value = 100
